$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.415.80'
$ws.Range('E2').Value = '  -2.89%  '
$ws.Range('D3').Value = '1.744.32'
$ws.Range('E3').Value = '  -3.50%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'321.60"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.24%  '
$ws.Range('D6').Value = "'1.001"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = "'0.4217"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -8.23%  '
$ws.Range('D8').Value = "'0.3583"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.66%  '
$ws.Range('D9').Value = "'45.50"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('D10').Value = "'0.07415"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.00%  '
$ws.Range('D11').Value = "'1.112"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('D12').Value = "'1.002"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').Value = "'21.46"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.42%  '
$ws.Range('D14').Value = "'6.107"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.99%  '
$ws.Range('E15').Value = '  -3.65%  '
$ws.Range('D16').Value = '1.740.03'
$ws.Range('E16').Value = '  -3.52%  '
$ws.Range('E17').Value = '  -3.00%  '
$ws.Range('D18').Value = "'87.58"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.75%  '
$ws.Range('D19').Value = "'0.06067"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -9.69%  '
$ws.Range('D20').Value = "'1.001"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('E21').Value = '  -3.93%  '
$ws.Range('D22').Value = "'6.103"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.34%  '
$ws.Range('D23').Value = "'0.5239"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.21%  '
$ws.Range('D24').Value = '27.448.66'
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('D25').Value = "'11.47"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.81%  '
$ws.Range('D26').Value = "'2.340"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.11%  '
$ws.Range('D27').Value = "'20.40"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.18%  '
$ws.Range('D28').Value = "'2.376"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').Value = "'152.11"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('D30').Value = '1.937.24'
$ws.Range('E30').Value = '  -3.59%  '
$ws.Range('D31').Value = "'125.75"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.08%  '
$ws.Range('D32').Value = "'1.198"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.11%  '
$ws.Range('D33').Value = "'5.681"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.46%  '
$ws.Range('D34').Value = "'0.09125"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.98%  '
$ws.Range('D35').Value = "'3.626"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -10.14%  '
$ws.Range('D36').Value = "'12.68"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.15%  '
$ws.Range('D37').Value = "'0.02293"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.95%  '
$ws.Range('D38').Value = "'0.2142"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.74%  '
$ws.Range('D39').Value = "'5.086"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.66%  '
$ws.Range('D40').Value = "'0.06052"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.40%  '
$ws.Range('D41').Value = "'0.6388"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.10%  '
$ws.Range('D42').Value = "'1.189"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.27%  '
$ws.Range('D43').Value = "'1.419"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.08%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = "'1.000"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = "'7.911"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.52%  '
$ws.Range('D46').Value = "'13.74"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.51%  '
$ws.Range('E47').Value = '  -3.07%  '
$ws.Range('D48').Value = "'0.5854"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.80%  '
$ws.Range('D49').Value = "'125.24"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.92%  '
$ws.Range('D50').Value = "'1.947"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.32%  '
$ws.Range('E51').Value = '  -4.72%  '
